$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.420357823371887
$ws.Range("B1").Value = 1.840587139129639
$ws.Range("C1").Value = 2.057744026184082
$ws.Range("D1").Value = 2.3577880859375
$ws.Range("E1").Value = 2.871772527694702
